# Insert a new data row at row 265 of the single worksheet ("Hortaliza, Vega
# Central Mapocho de Santiago - Haba"). All existing rows from 265 downward
# shift down by one (old row 265 becomes row 266, ..., old row 374 becomes
# row 375), and the freshly inserted row 265 is populated with a new
# weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 265..374 down to 266..375, leaving row 265 empty and ready to
# receive the new record.
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with the new observation.
$ws.Cells.Item(265, 1).Value = 9
$ws.Cells.Item(265, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(265, 3).Value = "Metropolitana"
$ws.Cells.Item(265, 4).Value = 45146
$ws.Cells.Item(265, 5).Value = 13
$ws.Cells.Item(265, 6).Value = 100112026
$ws.Cells.Item(265, 7).Value = "Haba"
$ws.Cells.Item(265, 8).Value = "Sin especificar"
$ws.Cells.Item(265, 9).Value = "Primera"
$ws.Cells.Item(265, 10).Value = 52
$ws.Cells.Item(265, 11).Value = 13000
$ws.Cells.Item(265, 12).Value = 15000
$ws.Cells.Item(265, 13).Value = 14000
$ws.Cells.Item(265, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(265, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(265, 16).Value = 560
$ws.Cells.Item(265, 17).Value = 25
$ws.Cells.Item(265, 18).Value = "Hortaliza"
